# Adds two new columns, I (header "I0") and J (header "IF"), to Sheet1,
# populating the header row and the 52 data rows (rows 2-53) with the
# values from the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the bold / bordered style used by the rest of the header row (e.g. H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data for rows 2..53 -> columns I (index 9) and J (index 10)
$ijData = @(
    @(5,7),
    @(8,8),
    @(7,8),
    @(6,6),
    @(6,7),
    @(6,6),
    @(2,5),
    @(8,8),
    @(1,3),
    @(9,9),
    @(3,4),
    @(7,8),
    @(13,13),
    @(5,5),
    @(4,5),
    @(6,7),
    @(5,6),
    @(5,7),
    @(7,7),
    @(6,8),
    @(3,4),
    @(7,7),
    @(3,7),
    @(9,9),
    @(7,7),
    @(5,5),
    @(4,6),
    @(6,8),
    @(6,8),
    @(5,6),
    @(8,8),
    @(7,8),
    @(8,9),
    @(7,8),
    @(7,8),
    @(5,6),
    @(8,9),
    @(8,9),
    @(8,9),
    @(7,8),
    @(8,8),
    @(7,8),
    @(2,6),
    @(5,5),
    @(6,7),
    @(9,9),
    @(4,5),
    @(8,9),
    @(8,8),
    @(9,9),
    @(7,8),
    @(5,6)
)

$row = 2
foreach ($pair in $ijData) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
